# Auto-generated: apply crypto price/volume updates from diff
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$updates = @(
    @{ Cell = 'D2'; Value = '60.785.80' },
    @{ Cell = 'E2'; Value = '  -1.15%  ' },
    @{ Cell = 'D3'; Value = '3.371.09' },
    @{ Cell = 'E3'; Value = '  -0.46%  ' },
    @{ Cell = 'D4'; Value = '1.00' },
    @{ Cell = 'E4'; Value = '  -0.05%  ' },
    @{ Cell = 'D5'; Value = '569.12' },
    @{ Cell = 'E5'; Value = '  -1.19%  ' },
    @{ Cell = 'D6'; Value = '135.87' },
    @{ Cell = 'E6'; Value = '  -0.74%  ' },
    @{ Cell = 'E7'; Value = '  +0.05%  ' },
    @{ Cell = 'D8'; Value = '3.369.44' },
    @{ Cell = 'E8'; Value = '  -0.51%  ' },
    @{ Cell = 'E9'; Value = '  -1.04%  ' },
    @{ Cell = 'D10'; Value = '7.59' },
    @{ Cell = 'E10'; Value = '  +1.14%  ' },
    @{ Cell = 'E11'; Value = '  -3.28%  ' },
    @{ Cell = 'E12'; Value = '  -2.72%  ' },
    @{ Cell = 'D13'; Value = '3.947.79' },
    @{ Cell = 'E13'; Value = '  -0.47%  ' },
    @{ Cell = 'E14'; Value = '  -0.77%  ' },
    @{ Cell = 'D15'; Value = '25.97' },
    @{ Cell = 'E15'; Value = '  +0.69%  ' },
    @{ Cell = 'D16'; Value = '3.373.39' },
    @{ Cell = 'E16'; Value = '  -0.31%  ' },
    @{ Cell = 'D17'; Value = '0.0000169' },
    @{ Cell = 'E17'; Value = '  -3.74%  ' },
    @{ Cell = 'D18'; Value = '60.839.84' },
    @{ Cell = 'E18'; Value = '  -1.22%  ' },
    @{ Cell = 'D19'; Value = '5.81' },
    @{ Cell = 'E19'; Value = '  -1.00%  ' },
    @{ Cell = 'D20'; Value = '13.71' },
    @{ Cell = 'E20'; Value = '  -3.32%  ' },
    @{ Cell = 'D21'; Value = '9.20' },
    @{ Cell = 'E21'; Value = '  -2.36%  ' },
    @{ Cell = 'D22'; Value = '371.15' },
    @{ Cell = 'E22'; Value = '  -1.33%  ' },
    @{ Cell = 'D23'; Value = '3.510.25' },
    @{ Cell = 'E23'; Value = '  -0.61%  ' },
    @{ Cell = 'E24'; Value = '  -1.81%  ' },
    @{ Cell = 'E25'; Value = '  +0.01%  ' },
    @{ Cell = 'D26'; Value = '70.67' },
    @{ Cell = 'E26'; Value = '  -0.80%  ' },
    @{ Cell = 'D27'; Value = '0.0000122' },
    @{ Cell = 'E27'; Value = '  -3.49%  ' },
    @{ Cell = 'E28'; Value = '  +8.68%  ' },
    @{ Cell = 'D29'; Value = '1.56' },
    @{ Cell = 'E29'; Value = '  -7.75%  ' },
    @{ Cell = 'D30'; Value = '1.00' },
    @{ Cell = 'E30'; Value = '  -0.09%  ' },
    @{ Cell = 'D31'; Value = '7.32' },
    @{ Cell = 'E31'; Value = '  -2.57%  ' },
    @{ Cell = 'D32'; Value = '8.00' },
    @{ Cell = 'E32'; Value = '  -2.78%  ' },
    @{ Cell = 'E33'; Value = '  -2.57%  ' },
    @{ Cell = 'E34'; Value = '  -0.05%  ' },
    @{ Cell = 'D35'; Value = '23.23' },
    @{ Cell = 'E35'; Value = '  -0.90%  ' },
    @{ Cell = 'D36'; Value = '5.08' },
    @{ Cell = 'E36'; Value = '  -4.37%  ' },
    @{ Cell = 'E37'; Value = '  -1.06%  ' },
    @{ Cell = 'D38'; Value = '6.75' },
    @{ Cell = 'E38'; Value = '  -1.12%  ' },
    @{ Cell = 'D39'; Value = '164.53' },
    @{ Cell = 'E39'; Value = '  -0.53%  ' },
    @{ Cell = 'D40'; Value = '0.0756' },
    @{ Cell = 'E40'; Value = '  -3.04%  ' },
    @{ Cell = 'E41'; Value = '  +0.02%  ' },
    @{ Cell = 'D42'; Value = '1.73' },
    @{ Cell = 'E42'; Value = '  +0.72%  ' },
    @{ Cell = 'D43'; Value = '0.770' },
    @{ Cell = 'E43'; Value = '  -1.36%  ' },
    @{ Cell = 'D44'; Value = '41.87' },
    @{ Cell = 'E44'; Value = '  +1.04%  ' },
    @{ Cell = 'D45'; Value = '25.04' },
    @{ Cell = 'E45'; Value = '  +0.09%  ' },
    @{ Cell = 'E46'; Value = '  -2.13%  ' },
    @{ Cell = 'D47'; Value = '1.16' },
    @{ Cell = 'E47'; Value = '  -5.85%  ' },
    @{ Cell = 'D48'; Value = '2.505.22' },
    @{ Cell = 'E48'; Value = '  +7.34%  ' },
    @{ Cell = 'D49'; Value = '23.40' },
    @{ Cell = 'E49'; Value = '  +3.59%  ' },
    @{ Cell = 'E50'; Value = '  -1.59%  ' },
    @{ Cell = 'D51'; Value = '2.38' },
    @{ Cell = 'E51'; Value = '  -0.15%  ' }
)

foreach ($u in $updates) {
    $rng = $ws.Range($u.Cell)
    $rng.NumberFormat = "@"
    $rng.Value = $u.Value
    $rng.Style = "Normal"
}

